$d = $word.ActiveDocument

# --- Step 1: remove the "Meta description" paragraph that follows the title ---
$titlePara = $d.Paragraphs.Item(1)
$metaPara = $titlePara.Next()
[void]$metaPara.Range.Delete()

# --- Step 2: insert a new bold "Play Book of Shadows..." paragraph just before
#     the final (image-prompt) paragraph ---
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
[void]$lastPara.Range.InsertParagraphBefore()
$newCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newCount - 1)
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:r/>' + `
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Play Book of Shadows for Free ' + [char]0x2013 + ' Review of Gameplay &amp; Bonus Features</w:t></w:r>' + `
    '</w:p>'
[void]$newPara.Range.InsertXML($titleXml)

# --- Step 3: replace the text of the final paragraph (formerly the image-generation
#     prompt) with the meta description text, keeping its italic run formatting ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$oldText = "Create a feature image for Book of Shadows, a horror-themed slot game, that is in cartoon style and features a happy Maya warrior with glasses. The image should showcase the Maya warrior holding the Book of Shadows with a confident and mischievous expression on his face, ready to tackle the horrors in the game. The background should depict a dark, eerie forest with moonlight casting a shadowy glow. The Maya warrior should be depicted wearing traditional Maya clothing, including a headdress adorned with feathers. The glasses should be modern, adding a playful touch to the image. The overall style should be a mix of ancient and modern, representing the theme of the game. The image should be colorful, bold, and attention-grabbing to entice players to try the game."
$newText = "Read our review of Book of Shadows and play for free today. Learn about the gameplay mechanics and bonus features, including the Free Spins and Nolimit bonus."
[void]$lastPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
